$d = $word.ActiveDocument
$d.Content.Find.Execute("Large Class", $true, $false, $false, $false, $false, $true, 1, $false, "Large Methods", 2)
